$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 280.1474913333333
$ws.Range("H2").Value = 840.4424739999999
$ws.Range("I2").Value = 0.298539071964842
$ws.Range("J2").Value = 0.298539071964842
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.709791333333333
$ws.Range("N2").Value = 11.129374
$ws.Range("O2").Value = 0.4283284425582907
$ws.Range("P2").Value = 0.4283284425582907
$ws.Range("Q2").Value = 1039.288735403475
$ws.Range("R2").Value = 9353.598618631275
$ws.Range("S2").Value = 0.1278727757374982
$ws.Range("T2").Value = 0.1278727757374982
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 280.1474913333333
$ws.Range("H3").Value = 840.4424739999999
$ws.Range("I3").Value = 0.298539071964842
$ws.Range("J3").Value = 0.298539071964842
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.599001333333334
$ws.Range("N3").Value = 10.797004
$ws.Range("O3").Value = 0.4155367505499981
$ws.Range("P3").Value = 0.4155367505499982
$ws.Range("Q3").Value = 1008.251194838655
$ws.Range("R3").Value = 9074.260753547896
$ws.Range("S3").Value = 0.1240539558764825
$ws.Range("T3").Value = 0.1240539558764825
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 280.1474913333333
$ws.Range("H4").Value = 840.4424739999999
$ws.Range("I4").Value = 0.298539071964842
$ws.Range("J4").Value = 0.298539071964842
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.352297666666667
$ws.Range("N4").Value = 4.056893
$ws.Range("O4").Value = 0.1561348068917112
$ws.Range("P4").Value = 0.1561348068917112
$ws.Range("Q4").Value = 378.8427988525869
$ws.Range("R4").Value = 3409.585189673282
$ws.Range("S4").Value = 0.04661234035086127
$ws.Range("T4").Value = 0.04661234035086127
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 136.0425923333333
$ws.Range("H5").Value = 408.127777
$ws.Range("I5").Value = 0.1449737389029841
$ws.Range("J5").Value = 0.1449737389029841
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 3.709791333333333
$ws.Range("N5").Value = 11.129374
$ws.Range("O5").Value = 0.4283284425582907
$ws.Range("P5").Value = 0.4283284425582907
$ws.Range("Q5").Value = 504.6896300023997
$ws.Range("R5").Value = 4542.206670021598
$ws.Range("S5").Value = 0.06209637579616744
$ws.Range("T5").Value = 0.06209637579616745
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 136.0425923333333
$ws.Range("H6").Value = 408.127777
$ws.Range("I6").Value = 0.1449737389029841
$ws.Range("J6").Value = 0.1449737389029841
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.599001333333334
$ws.Range("N6").Value = 10.797004
$ws.Range("O6").Value = 0.4155367505499981
$ws.Range("P6").Value = 0.4155367505499982
$ws.Range("Q6").Value = 489.6174711977898
$ws.Range("R6").Value = 4406.557240780108
$ws.Range("S6").Value = 0.06024191637882985
$ws.Range("T6").Value = 0.06024191637882986
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 136.0425923333333
$ws.Range("H7").Value = 408.127777
$ws.Range("I7").Value = 0.1449737389029841
$ws.Range("J7").Value = 0.1449737389029841
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.352297666666667
$ws.Range("N7").Value = 4.056893
$ws.Range("O7").Value = 0.1561348068917112
$ws.Range("P7").Value = 0.1561348068917112
$ws.Range("Q7").Value = 183.9700801796512
$ws.Range("R7").Value = 1655.730721616861
$ws.Range("S7").Value = 0.02263544672798678
$ws.Range("T7").Value = 0.02263544672798678
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 206.9361446666667
$ws.Range("H8").Value = 620.808434
$ws.Range("I8").Value = 0.2205214270909241
$ws.Range("J8").Value = 0.2205214270909241
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 3.709791333333333
$ws.Range("N8").Value = 11.129374
$ws.Range("O8").Value = 0.4283284425582907
$ws.Range("P8").Value = 0.4283284425582907
$ws.Range("Q8").Value = 767.689916037813
$ws.Range("R8").Value = 6909.209244340316
$ws.Range("S8").Value = 0.09445559941658717
$ws.Range("T8").Value = 0.09445559941658717
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 206.9361446666667
$ws.Range("H9").Value = 620.808434
$ws.Range("I9").Value = 0.2205214270909241
$ws.Range("J9").Value = 0.2205214270909241
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.599001333333334
$ws.Range("N9").Value = 10.797004
$ws.Range("O9").Value = 0.4155367505499981
$ws.Range("P9").Value = 0.4155367505499982
$ws.Range("Q9").Value = 744.763460570193
$ws.Range("R9").Value = 6702.871145131737
$ws.Range("S9").Value = 0.09163475724001093
$ws.Range("T9").Value = 0.09163475724001093
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 206.9361446666667
$ws.Range("H10").Value = 620.808434
$ws.Range("I10").Value = 0.2205214270909241
$ws.Range("J10").Value = 0.2205214270909241
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.352297666666667
$ws.Range("N10").Value = 4.056893
$ws.Range("O10").Value = 0.1561348068917112
$ws.Range("P10").Value = 0.1561348068917112
$ws.Range("Q10").Value = 279.8392655817291
$ws.Range("R10").Value = 2518.553390235562
$ws.Range("S10").Value = 0.034431070434326
$ws.Range("T10").Value = 0.034431070434326
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 305.6466113333333
$ws.Range("H11").Value = 916.939834
$ws.Range("I11").Value = 0.325712199892882
$ws.Range("J11").Value = 0.325712199892882
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 3.709791333333333
$ws.Range("N11").Value = 11.129374
$ws.Range("O11").Value = 0.4283284425582907
$ws.Range("P11").Value = 0.4283284425582907
$ws.Range("Q11").Value = 1133.885149787102
$ws.Range("R11").Value = 10204.96634808392
$ws.Range("S11").Value = 0.1395117993023528
$ws.Range("T11").Value = 0.1395117993023528
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 305.6466113333333
$ws.Range("H12").Value = 916.939834
$ws.Range("I12").Value = 0.325712199892882
$ws.Range("J12").Value = 0.325712199892882
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 3.599001333333334
$ws.Range("N12").Value = 10.797004
$ws.Range("O12").Value = 0.4155367505499981
$ws.Range("P12").Value = 0.4155367505499982
$ws.Range("Q12").Value = 1100.022561717482
$ws.Range("R12").Value = 9900.203055457338
$ws.Range("S12").Value = 0.1353453891579796
$ws.Range("T12").Value = 0.1353453891579796
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 305.6466113333333
$ws.Range("H13").Value = 916.939834
$ws.Range("I13").Value = 0.325712199892882
$ws.Range("J13").Value = 0.325712199892882
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 1.352297666666667
$ws.Range("N13").Value = 4.056893
$ws.Range("O13").Value = 0.1561348068917112
$ws.Range("P13").Value = 0.1561348068917112
$ws.Range("Q13").Value = 413.3251993306402
$ws.Range("R13").Value = 3719.926793975762
$ws.Range("S13").Value = 0.05085501143254956
$ws.Range("T13").Value = 0.05085501143254956
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 9.621888666666665
$ws.Range("H14").Value = 28.865666
$ws.Range("I14").Value = 0.01025356214836792
$ws.Range("J14").Value = 0.01025356214836792
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 3.709791333333333
$ws.Range("N14").Value = 11.129374
$ws.Range("O14").Value = 0.4283284425582907
$ws.Range("P14").Value = 0.4283284425582907
$ws.Range("Q14").Value = 35.69519918589822
$ws.Range("R14").Value = 321.256792673084
$ws.Range("S14").Value = 0.004391892305685072
$ws.Range("T14").Value = 0.004391892305685073
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 9.621888666666665
$ws.Range("H15").Value = 28.865666
$ws.Range("I15").Value = 0.01025356214836792
$ws.Range("J15").Value = 0.01025356214836792
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 3.599001333333334
$ws.Range("N15").Value = 10.797004
$ws.Range("O15").Value = 0.4155367505499981
$ws.Range("P15").Value = 0.4155367505499982
$ws.Range("Q15").Value = 34.62919014051822
$ws.Range("R15").Value = 311.662711264664
$ws.Range("S15").Value = 0.004260731896695263
$ws.Range("T15").Value = 0.004260731896695264
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 9.621888666666665
$ws.Range("H16").Value = 28.865666
$ws.Range("I16").Value = 0.01025356214836792
$ws.Range("J16").Value = 0.01025356214836792
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 1.352297666666667
$ws.Range("N16").Value = 4.056893
$ws.Range("O16").Value = 0.1561348068917112
$ws.Range("P16").Value = 0.1561348068917112
$ws.Range("Q16").Value = 13.01165759285978
$ws.Range("R16").Value = 117.104918335738
$ws.Range("S16").Value = 0.001600937945987585
$ws.Range("T16").Value = 0.001600937945987585
